# Actualizacion automatica del mapa (2025-12-12 15:25:36)
#
# The source "Optical_Power" log sheet had its oldest pending-claim row
# (old row 89 - Caso 7853 / ACOSTA, MARIANO 2769) resolved/removed, which
# shifts every following record up by one row. One corrected OT value is
# applied (new row 92, column E), and two freshly logged claims are
# appended at the bottom (new rows 94 and 95).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that are stored as TEXT in this sheet (case numbers, OT
# numbers, comuna codes, mm/dd/yyyy date strings, free-text notes, ...).
# Columns I, M, N are genuine numbers.
$textCols = @("A","B","C","D","E","F","G","H","J","K","L","O","P","Q","R")

function Looks-NumericOrDate($value) {
    if ($value -eq $null) { return $false }
    $s = [string]$value
    if ($s.Trim().Length -eq 0) { return $false }
    if ($s -match '^\d{1,2}/\d{1,2}/\d{4}$') { return $true }
    if ($s.Trim() -match '^-?\d+(\.\d+)?$') { return $true }
    return $false
}

function Set-TextCell($row, $col, $value) {
    # Only force the cell to text format when the literal value would
    # otherwise be auto-converted by Excel into a number/date (which
    # would corrupt leading zeros / trailing spaces / string-ness).
    $cell = $ws.Range($col + $row)
    if (Looks-NumericOrDate $value) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

function Set-RowData($row, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $val = $values[$i]
        if ($textCols -contains $col) {
            Set-TextCell $row $col $val
        } else {
            $ws.Range($col + $row).Value = $val
        }
    }
}

# 1) The old row 89 (Caso 7853, ACOSTA, MARIANO 2769) is gone; every row
#    below it moves up one slot. Deleting the row performs exactly that
#    shift while preserving each remaining record's own data untouched.
$ws.Rows.Item(89).Delete()

# 2) After the shift, the record that is now row 92 (Caso 7892, Juncal
#    1642) gets its OT corrected from "Pendiente ADM" to the real OT.
Set-TextCell 92 "E" "01749376 "

# 3) Two new claims logged at the bottom of the sheet.
Set-RowData 94 @("7879 ", "12/9/2025", "CORVALAN 996", "9", "811131632", "Optical Power", "Pendiente", "base corroida", 1, "Cambio", "Sin equipos", "Pasante", -58.498487, -34.647524, "Devoto", "Capital Norte", "PAV-?", "Fuera de Poligono OVL")

Set-RowData 95 @("S01204545", "12/10/2025", "Medina 420", "10", "811131640", "Optical Power", "Pendiente", "inclinado", 1, "Cambio", "Sin equipos", "Pasante", -58.48802, -34.641075, "Devoto", "Capital Norte", "PCH-S", "ARATO-25058.PO.2PCH")
